# Applies the "2020 column" update to the ИПЦ worksheet:
#  - adds column Q (year 2020) data for rows 4-14
#  - moves the active-cell selection to N14 (matches the author's saved
#    cursor position at commit time)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from column P (the previous "latest year" column) into
# column Q for the header + every data row, then fill in the 2020 values.
$ws.Range("P4:P14").Copy() | Out-Null
$ws.Range("Q4:Q14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("Q4").Value = 2020

$ws.Range("Q5").Value = 109.7221295941265
$ws.Range("Q6").Value = 108.44905375816947
$ws.Range("Q7").Value = 109.90982951756889
$ws.Range("Q8").Value = 108.40606487500015
$ws.Range("Q9").Value = 109.40161876466024
$ws.Range("Q10").Value = 107.71155656686271
$ws.Range("Q11").Value = 111.78921596090774
$ws.Range("Q12").Value = 111.39254046803097
$ws.Range("Q13").Value = 110.44919152842827
$ws.Range("Q14").Value = 106.89826464456031

# Restore the saved cursor / selection position recorded in the workbook.
$ws.Range("N14").Select() | Out-Null
